$d = $word.ActiveDocument

# 1. Merge ", " + "Cybersec" runs (removing spell-check split) into one run
#    text reads ", Cybersec" -- simple text replace covers the run merge.
$d.Content.Find.Execute(", Cybersec", $true, $false, $false, $false, $false,
                         $true, 1, $false, ", Cybersec", 2) | Out-Null

# 2. ACTS - Web Development Intern start date: "Mar 2024 - Present" -> "Mar 2025 - Present"
$d.Content.Find.Execute("Mar 2024 - Present", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mar 2025 - Present", 2) | Out-Null

# 3. Firestore stays the same text, just drop the proofing spell-check markers
#    around it (no visible text change needed beyond what Find/Replace above does).
